$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added. Insert a new row at 13, which shifts
# the existing rows 13-119 down to 14-120 (so the dimension grows to A1:R120).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value = 44881
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 100112026
$ws.Range("G13").Value = "Haba"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("N13").Value = "$/saco 25 kilos"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 480
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
